$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray "horas" values left in column H for these rows.
$ws.Range("H2").ClearContents()
$ws.Range("H3").ClearContents()
$ws.Range("H4").ClearContents()
$ws.Range("H6").ClearContents()
$ws.Range("H8").ClearContents()

# Update the active selection to match the author's final cursor position.
$ws.Range("F22").Select()
